$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated BOM table (refreshed from the linked TS3A5017-Mixer.csv data source):
# silkscreen / ordering rework -> new designators, footprints and comments.
$data = New-Object 'object[,]' 12,3

$data[0,0] = "Comment";        $data[0,1] = "Designator";                                    $data[0,2] = "Footprint"
$data[1,0] = ".1uf";           $data[1,1] = "C4,C5,C10,C9,C3,C12,C6,C8,C13,C11,C2,C7";        $data[1,2] = "C_0201_0603Metric"
$data[2,0] = ".1uf";           $data[2,1] = "C1";                                             $data[2,2] = "C_0603_1608Metric"
$data[3,0] = "TS3A5017RGY";    $data[3,1] = "U1";                                             $data[3,2] = "Texas_RGY_R-PVQFN-N16_EP2.05x2.55mm"
$data[4,0] = "CLK1";           $data[4,1] = "J3";                                             $data[4,2] = "SMA_Amphenol_901-143_Horizontal"
$data[5,0] = "Conn_01x02_Pin"; $data[5,1] = "J4";                                             $data[5,2] = "PinHeader_1x02_P2.54mm_Vertical"
$data[6,0] = "1uf";            $data[6,1] = "C14,C15";                                        $data[6,2] = "C_0201_0603Metric"
$data[7,0] = "RF-In";          $data[7,1] = "J1";                                             $data[7,2] = "SMA_Amphenol_901-143_Horizontal"
$data[8,0] = "CLK0";           $data[8,1] = "J2";                                             $data[8,2] = "SMA_Amphenol_901-143_Horizontal"
$data[9,0] = "10k";            $data[9,1] = "R1,R2";                                          $data[9,2] = "R_0201_0603Metric"
$data[10,0] = "LMP7715MF";     $data[10,1] = "U2,U3";                                         $data[10,2] = "SOT-23-5"
$data[11,0] = "PJ-320D-A";     $data[11,1] = "J5";                                            $data[11,2] = "HRO_PJ-320D-A"

# Clear the old table body first (old sheet had only 8 rows, new one has 12).
$ws.Range("A1:C12").ClearContents()

$ws.Range("A1:C12").Value = $data

# The named range driving the BOM table shrank from 8 rows of data to 5.
$name = $wb.Names.Item("TS3A5017_Mixer")
$name.RefersTo = "=Sheet1!`$A`$1:`$E`$5"

# Reflect the author's last cursor position when they saved the file.
$ws.Range("G9").Select()
